$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace(
    "1000 Bs = 7.02 = 27404.29 pesos",
    "1000 Bs = 6.95 = 27198.83 pesos"
)
$text = $text.Replace(
    "27404.29 pesos = 6.99 = 978.32 Bs",
    "27198.83 pesos = 6.91 = 966.95 Bs"
)
$cell.Value = $text

# --- Sheet "tasas": update the scraped rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 143.94
$wsTasas.Range("O10").Value = 3915
$wsTasas.Range("N12").Value = 3937.99
$wsTasas.Range("O12").Value = 140
